# Actualización automática 2025-07-01 14:25:07
$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M2").Value = 66.2
$wsGrupo.Range("M9").Value = "1 de 7"

$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F2").Value = 66.2
$wsMensual.Range("F9").Value = 66.2
